$wb = $excel.ActiveWorkbook

# ---- Sheet: Productdata ----
$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 15.496
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 5.952000000000001
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 4.0128
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 1.672
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 3.0096
$ws.Range("C7").Value = 899
$ws.Range("E7").Value = 13.156
$ws.Range("C8").Value = 138
$ws.Range("E8").Value = 4.572
$ws.Range("C9").Value = 232
$ws.Range("E9").Value = 3.72

# ---- Sheet: Capacity ----
$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 8940
$ws.Range("B3").Value = 6000
$ws.Range("B4").Value = 12540
$ws.Range("B5").Value = 12540
$ws.Range("B6").Value = 4180
$ws.Range("B7").Value = 10120
$ws.Range("B8").Value = 1350
$ws.Range("B9").Value = 1500

# ---- Sheet: ProcessingTime ----
$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("B2").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("H8").Value = 3
$ws.Range("I9").Value = 2
